# Rename 'terminations' to 'targets'
#
# Also reproduces the view-state side effects captured in the same save
# (the renamed sheet becomes the active/selected tab with a new selection,
# and the 'flow' sheet's frozen-pane view is scrolled further down).

$wb = $excel.ActiveWorkbook

# --- Rename the 'terminations' worksheet to 'targets' -----------------
$wsTargets = $wb.Worksheets.Item("terminations")
$wsTargets.Name = "targets"

# --- Scroll the 'flow' sheet's frozen pane down (best effort) ---------
$wsFlow = $wb.Worksheets.Item("flow")
$wsFlow.Activate()
$winFlow = $excel.ActiveWindow
$winFlow.ScrollRow = 200
$winFlow.ScrollColumn = 1

# --- Make 'targets' the active sheet, with its own new selection ------
$wsTargets.Activate()
$wsTargets.Range("E29").Select() | Out-Null
